# Appends daily rows 329..343 (dates 2021-07-26 .. 2021-08-09, serials
# 44403..44417) to the bottom of the existing data table, matching the
# layout/style of the preceding rows (column A keeps the date style,
# columns B/C/D are plain numeric zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 328
$startDate = 44403
$numNewRows = 15

for ($i = 0; $i -lt $numNewRows; $i++) {
    $row = $lastRow + 1 + $i
    $date = $startDate + $i

    # Copy the full row above so formatting (date style on col A, default
    # style on B:D) carries forward, then overwrite with the new values.
    $srcRow = $row - 1
    $src = $ws.Range("A" + $srcRow + ":D" + $srcRow)
    $dst = $ws.Range("A" + $row + ":D" + $row)
    $src.Copy($dst)

    $ws.Range("A" + $row).Value = $date
    $ws.Range("B" + $row).Value = 0
    $ws.Range("C" + $row).Value = 0
    $ws.Range("D" + $row).Value = 0
}
